$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, "H").Value = 153.28572
$ws.Cells.Item(4, "I").Value = 153.28572
$ws.Cells.Item(4, "K").Value = 153.28572
$ws.Cells.Item(4, "M").Value = -39.28572

$ws.Cells.Item(9, "H").Value = 191.42857
$ws.Cells.Item(9, "J").Value = 224
$ws.Cells.Item(9, "L").Value = 224
$ws.Cells.Item(9, "N").Value = -562

$ws.Cells.Item(29, "H").Value = 1124.6875
$ws.Cells.Item(29, "J").Value = 899.6429000000001
$ws.Cells.Item(29, "L").Value = 2698.9287
$ws.Cells.Item(29, "N").Value = -3260.9287

$ws.Cells.Item(33, "H").Value = 140.16667
$ws.Cells.Item(33, "I").Value = 156.42857
$ws.Cells.Item(33, "J").Value = 83.25
$ws.Cells.Item(33, "K").Value = 156.42857
$ws.Cells.Item(33, "L").Value = 83.25
$ws.Cells.Item(33, "M").Value = 72.57142999999999
$ws.Cells.Item(33, "N").Value = -541.25

$ws.Cells.Item(42, "H").Value = 589.9091
$ws.Cells.Item(42, "I").Value = 20.25
$ws.Cells.Item(42, "K").Value = 60.75
$ws.Cells.Item(42, "M").Value = 169.25

$ws.Cells.Item(74, "H").Value = 4450
$ws.Cells.Item(74, "I").Value = 3500
$ws.Cells.Item(74, "J").Value = 4857.143
$ws.Cells.Item(74, "K").Value = 3500
$ws.Cells.Item(74, "L").Value = 4857.143
$ws.Cells.Item(74, "M").Value = -2564
$ws.Cells.Item(74, "N").Value = -6729.143

$ws.Cells.Item(77, "H").Value = 4450
$ws.Cells.Item(77, "I").Value = 3500
$ws.Cells.Item(77, "J").Value = 4857.143
$ws.Cells.Item(77, "K").Value = 17500
$ws.Cells.Item(77, "L").Value = 24285.715
$ws.Cells.Item(77, "M").Value = -12820
$ws.Cells.Item(77, "N").Value = -33645.715

$ws.Cells.Item(86, "H").Value = 6455573
$ws.Cells.Item(86, "I").Value = 16132032
$ws.Cells.Item(86, "J").Value = 4600
$ws.Cells.Item(86, "K").Value = 16132032
$ws.Cells.Item(86, "L").Value = 4600
$ws.Cells.Item(86, "M").Value = -16130909
$ws.Cells.Item(86, "N").Value = -6846

$ws.Cells.Item(88, "H").Value = 16678074
$ws.Cells.Item(88, "J").Value = 13431.889
$ws.Cells.Item(88, "L").Value = 13431.889
$ws.Cells.Item(88, "N").Value = -14243.889

$ws.Cells.Item(89, "H").Value = 6455573
$ws.Cells.Item(89, "I").Value = 16132032
$ws.Cells.Item(89, "J").Value = 4600
$ws.Cells.Item(89, "K").Value = 80660160
$ws.Cells.Item(89, "L").Value = 23000
$ws.Cells.Item(89, "M").Value = -80654544
$ws.Cells.Item(89, "N").Value = -34232

$ws.Cells.Item(91, "H").Value = 16678074
$ws.Cells.Item(91, "J").Value = 13431.889
$ws.Cells.Item(91, "L").Value = 13431.889
$ws.Cells.Item(91, "N").Value = -16239.889

$ws.Cells.Item(111, "H").Value = 4693.5
$ws.Cells.Item(111, "I").Value = 4924.6665
$ws.Cells.Item(111, "K").Value = 14773.9995
$ws.Cells.Item(111, "M").Value = -11706.9995

$ws.Cells.Item(113, "H").Value = 3403.8823
$ws.Cells.Item(113, "I").Value = 3718.3125
$ws.Cells.Item(113, "K").Value = 3718.3125
$ws.Cells.Item(113, "M").Value = -464.3125

$ws.Cells.Item(136, "H").Value = 2092395.8
$ws.Cells.Item(136, "J").Value = 2092395.8
$ws.Cells.Item(136, "L").Value = 2092395.8
$ws.Cells.Item(136, "N").Value = -2102595.8

$ws.Cells.Item(137, "H").Value = 4448.4
$ws.Cells.Item(137, "J").Value = 5612.4287
$ws.Cells.Item(137, "L").Value = 16837.2861
$ws.Cells.Item(137, "N").Value = -21937.2861

$ws.Cells.Item(138, "H").Value = 8082.0625
$ws.Cells.Item(138, "J").Value = 8399.200000000001
$ws.Cells.Item(138, "L").Value = 25197.6
$ws.Cells.Item(138, "N").Value = -35477.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, "H").Value = 2209.2727
$ws.Cells.Item(45, "I").Value = 2208
$ws.Cells.Item(45, "K").Value = 2208
$ws.Cells.Item(45, "M").Value = -1831

$ws.Cells.Item(61, "H").Value = 12420.556
$ws.Cells.Item(61, "I").Value = 9909.409
$ws.Cells.Item(61, "J").Value = 16366.643
$ws.Cells.Item(61, "K").Value = 9909.409
$ws.Cells.Item(61, "L").Value = 16366.643
$ws.Cells.Item(61, "M").Value = -9697.409
$ws.Cells.Item(61, "N").Value = -16790.643

$ws.Cells.Item(74, "H").Value = 1469.7646
$ws.Cells.Item(74, "I").Value = 769.8
$ws.Cells.Item(74, "J").Value = 1761.4166
$ws.Cells.Item(74, "K").Value = 769.8
$ws.Cells.Item(74, "L").Value = 1761.4166
$ws.Cells.Item(74, "M").Value = 104.2
$ws.Cells.Item(74, "N").Value = -3509.4166

$ws.Cells.Item(77, "H").Value = 1469.7646
$ws.Cells.Item(77, "I").Value = 769.8
$ws.Cells.Item(77, "J").Value = 1761.4166
$ws.Cells.Item(77, "K").Value = 3849
$ws.Cells.Item(77, "L").Value = 8807.083000000001
$ws.Cells.Item(77, "M").Value = 519
$ws.Cells.Item(77, "N").Value = -17543.083

$ws.Cells.Item(122, "H").Value = 3768.3333
$ws.Cells.Item(122, "I").Value = 2614.6428
$ws.Cells.Item(122, "J").Value = 6075.7144
$ws.Cells.Item(122, "K").Value = 7843.928400000001
$ws.Cells.Item(122, "L").Value = 18227.1432
$ws.Cells.Item(122, "M").Value = -5393.928400000001
$ws.Cells.Item(122, "N").Value = -23127.1432

$ws.Cells.Item(132, "H").Value = 4367.3
$ws.Cells.Item(132, "I").Value = 1823.25
$ws.Cells.Item(132, "J").Value = 14543.5
$ws.Cells.Item(132, "K").Value = 5469.75
$ws.Cells.Item(132, "L").Value = 43630.5
$ws.Cells.Item(132, "M").Value = -2939.75
$ws.Cells.Item(132, "N").Value = -48690.5

$ws.Cells.Item(136, "H").Value = 12420.556
$ws.Cells.Item(136, "I").Value = 9909.409
$ws.Cells.Item(136, "J").Value = 16366.643
$ws.Cells.Item(136, "K").Value = 29728.227
$ws.Cells.Item(136, "L").Value = 49099.929
$ws.Cells.Item(136, "M").Value = -27178.227
$ws.Cells.Item(136, "N").Value = -54199.929

$ws.Cells.Item(140, "H").Value = 111065.5
$ws.Cells.Item(140, "J").Value = 111065.5
$ws.Cells.Item(140, "L").Value = 111065.5
$ws.Cells.Item(140, "N").Value = -121425.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, "H").Value = 4518
$ws.Cells.Item(86, "J").Value = 4777
$ws.Cells.Item(86, "L").Value = 4777
$ws.Cells.Item(86, "N").Value = -7023

$ws.Cells.Item(89, "H").Value = 4518
$ws.Cells.Item(89, "J").Value = 4777
$ws.Cells.Item(89, "L").Value = 23885
$ws.Cells.Item(89, "N").Value = -35117

$ws.Cells.Item(94, "H").Value = 857162.9
$ws.Cells.Item(94, "J").Value = 1041.1428
$ws.Cells.Item(94, "L").Value = 1041.1428
$ws.Cells.Item(94, "N").Value = -1943.1428

$ws.Cells.Item(99, "H").Value = 6945774
$ws.Cells.Item(99, "I").Value = 6945774
$ws.Cells.Item(99, "K").Value = 6945774
$ws.Cells.Item(99, "M").Value = -6944276

$ws.Cells.Item(105, "H").Value = 1975.5555
$ws.Cells.Item(105, "I").Value = 1771
$ws.Cells.Item(105, "J").Value = 2180.111
$ws.Cells.Item(105, "K").Value = 1771
$ws.Cells.Item(105, "L").Value = 2180.111
$ws.Cells.Item(105, "M").Value = -24
$ws.Cells.Item(105, "N").Value = -5674.111

$ws.Cells.Item(107, "H").Value = 1857.7778
$ws.Cells.Item(107, "I").Value = 3186.1667
$ws.Cells.Item(107, "J").Value = 1193.5834
$ws.Cells.Item(107, "K").Value = 3186.1667
$ws.Cells.Item(107, "L").Value = 1193.5834
$ws.Cells.Item(107, "M").Value = -1266.1667
$ws.Cells.Item(107, "N").Value = -5033.5834

$ws.Cells.Item(134, "H").Value = 3400.682
$ws.Cells.Item(134, "I").Value = 2926.1875
$ws.Cells.Item(134, "J").Value = 4666
$ws.Cells.Item(134, "K").Value = 8778.5625
$ws.Cells.Item(134, "L").Value = 13998
$ws.Cells.Item(134, "M").Value = -6243.5625
$ws.Cells.Item(134, "N").Value = -19068

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, "H").Value = 2849.611
$ws.Cells.Item(31, "I").Value = 945.0909
$ws.Cells.Item(31, "K").Value = 945.0909
$ws.Cells.Item(31, "M").Value = -650.0909

$ws.Cells.Item(34, "H").Value = 2849.611
$ws.Cells.Item(34, "I").Value = 945.0909
$ws.Cells.Item(34, "K").Value = 945.0909
$ws.Cells.Item(34, "M").Value = -743.0909

$ws.Cells.Item(75, "H").Value = 100000
$ws.Cells.Item(75, "J").Value = 100000
$ws.Cells.Item(75, "L").Value = 100000
$ws.Cells.Item(75, "N").Value = -101996

$ws.Cells.Item(78, "H").Value = 100000
$ws.Cells.Item(78, "J").Value = 100000
$ws.Cells.Item(78, "L").Value = 300000
$ws.Cells.Item(78, "N").Value = -309984

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(96, "H").Value = 12249.75
$ws.Cells.Item(96, "J").Value = 12999.667
$ws.Cells.Item(96, "L").Value = 38999.001
$ws.Cells.Item(96, "N").Value = -43117.001

$ws.Cells.Item(122, "H").Value = 492.9
$ws.Cells.Item(122, "J").Value = 491.2857
$ws.Cells.Item(122, "L").Value = 4421.571300000001
$ws.Cells.Item(122, "N").Value = -9321.5713

$ws.Cells.Item(132, "H").Value = 4928.2856
$ws.Cells.Item(132, "J").Value = 4300
$ws.Cells.Item(132, "L").Value = 38700
$ws.Cells.Item(132, "N").Value = -43760

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, "H").Value = 0
$ws.Cells.Item(98, "I").Value = 0
$ws.Cells.Item(98, "K").Value = 0
$ws.Cells.Item(98, "M").ClearContents()

$ws.Cells.Item(102, "H").Value = 10979.2
$ws.Cells.Item(102, "I").Value = 11724.083
$ws.Cells.Item(102, "K").Value = 11724.083
$ws.Cells.Item(102, "M").Value = -10102.083

$ws.Cells.Item(113, "H").Value = 10839.857
$ws.Cells.Item(113, "I").Value = 6826.3335
$ws.Cells.Item(113, "K").Value = 6826.3335
$ws.Cells.Item(113, "M").Value = -4656.3335

$ws.Cells.Item(126, "H").Value = 3962.6667
$ws.Cells.Item(126, "I").Value = 2018.091
$ws.Cells.Item(126, "J").Value = 5608.077
$ws.Cells.Item(126, "K").Value = 6054.272999999999
$ws.Cells.Item(126, "L").Value = 16824.231
$ws.Cells.Item(126, "M").Value = -3584.272999999999
$ws.Cells.Item(126, "N").Value = -21764.231

$ws.Cells.Item(132, "H").Value = 2388.4614
$ws.Cells.Item(132, "I").Value = 2547.2727
$ws.Cells.Item(132, "J").Value = 1515
$ws.Cells.Item(132, "K").Value = 7641.8181
$ws.Cells.Item(132, "L").Value = 4545
$ws.Cells.Item(132, "M").Value = -5111.8181
$ws.Cells.Item(132, "N").Value = -9605

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, "H").Value = 1833.3334
$ws.Cells.Item(93, "I").Value = 2000
$ws.Cells.Item(93, "K").Value = 2000
$ws.Cells.Item(93, "M").Value = -752

$ws.Cells.Item(122, "H").Value = 19889.6
$ws.Cells.Item(122, "I").Value = 5499.3335
$ws.Cells.Item(122, "K").Value = 16498.0005
$ws.Cells.Item(122, "M").Value = -14048.0005

$ws.Cells.Item(136, "H").Value = 7400
$ws.Cells.Item(136, "I").Value = 6100
$ws.Cells.Item(136, "K").Value = 18300
$ws.Cells.Item(136, "M").Value = -15750

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, "H").Value = 1159.3182
$ws.Cells.Item(107, "I").Value = 1287.6
$ws.Cells.Item(107, "J").Value = 884.4286
$ws.Cells.Item(107, "K").Value = 3862.8
$ws.Cells.Item(107, "L").Value = 2653.2858
$ws.Cells.Item(107, "M").Value = -1942.8
$ws.Cells.Item(107, "N").Value = -6493.2858

$ws.Cells.Item(132, "H").Value = 16132645
$ws.Cells.Item(132, "I").Value = 1286
$ws.Cells.Item(132, "J").Value = 166692000
$ws.Cells.Item(132, "K").Value = 3858
$ws.Cells.Item(132, "L").Value = 500076000
$ws.Cells.Item(132, "M").Value = -1328
$ws.Cells.Item(132, "N").Value = -500081060

$ws.Cells.Item(136, "H").Value = 7939.319
$ws.Cells.Item(136, "I").Value = 3843.6316
$ws.Cells.Item(136, "J").Value = 9020.125
$ws.Cells.Item(136, "K").Value = 11530.8948
$ws.Cells.Item(136, "L").Value = 27060.375
$ws.Cells.Item(136, "M").Value = -8980.8948
$ws.Cells.Item(136, "N").Value = -32160.375
